$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tab-free pipe-delimited data for the new arrival rows (row|A|B|C|D|E|F|G|H|I|J|L)
$data = @"
425|424.0|Monday, Jan 09|1:55 PM|LO282|London|(LHR)|LOT |B38M|(SP-LVB)|2:02 PM|0 hours, 7 minutes
426|425.0|Monday, Jan 09|1:55 PM|LX1348|Zurich|(ZRH)|Air Baltic |BCS3|(YL-ABL)|2:08 PM|0 hours, 13 minutes
427|426.0|Monday, Jan 09|2:00 PM|LH1348|Frankfurt|(FRA)|Lufthansa |A319|(D-AIBA)|2:40 PM|0 hours, 40 minutes
428|427.0|Monday, Jan 09|2:10 PM|LO536|Budapest|(BUD)|LOT (Star Alliance Livery) |E75S|(SP-LIO)|1:47 PM|0 hours, -23 minutes
429|428.0|Monday, Jan 09|2:10 PM|LO784|Riga|(RIX)|LOT |E75S|(SP-LIB)|2:10 PM|0 hours, 0 minutes
430|429.0|Monday, Jan 09|2:15 PM|LO528|Prague|(PRG)|LOT |E195|(SP-LNI)|1:53 PM|0 hours, -22 minutes
431|430.0|Monday, Jan 09|2:35 PM|LO156|Tel Aviv|(TLV)|LOT |B738|(SP-LWC)|2:31 PM|0 hours, -4 minutes
432|431.0|Monday, Jan 09|2:35 PM|W61576|Kutaisi|(KUT)|Wizz Air |A321|(HA-LTI)|2:25 PM|0 hours, -10 minutes
433|432.0|Monday, Jan 09|2:40 PM|LO786|Tallinn|(TLL)|LOT (Retro Livery) |E75S|(SP-LIM)|2:23 PM|0 hours, -17 minutes
434|433.0|Monday, Jan 09|2:50 PM|LO772|Vilnius|(VNO)|LOT |E75S|(SP-LIK)|2:49 PM|0 hours, -1 minutes
435|434.0|Monday, Jan 09|2:50 PM|W61560|Tel Aviv|(TLV)|Wizz Air |A21N|(HA-LVG)|2:45 PM|0 hours, -5 minutes
436|435.0|Monday, Jan 09|3:15 PM|LO616|Ljubljana|(LJU)|LOT |E75S|(SP-LIL)|2:47 PM|0 hours, -28 minutes
437|436.0|Monday, Jan 09|3:20 PM|LO642|Bucharest|(OTP)|LOT |E195|(SP-LND)|3:00 PM|0 hours, -20 minutes
438|437.0|Monday, Jan 09|3:24 PM|LO6592|Phu Quoc|(PQC)|LOT |B789|(SP-LSA)|3:26 PM|0 hours, 2 minutes
439|438.0|Monday, Jan 09|3:25 PM|BA842|London|(LHR)|British Airways |A320|(G-EUUA)|3:15 PM|0 hours, -10 minutes
440|439.0|Monday, Jan 09|3:30 PM|LO150|Cairo|(CAI)|LOT |B738|(SP-LWD)|3:37 PM|0 hours, 7 minutes
441|440.0|Monday, Jan 09|3:30 PM|LO514|Chisinau|(KIV)|LOT |B38M|(SP-LVA)|3:28 PM|0 hours, -2 minutes
442|441.0|Monday, Jan 09|3:35 PM|LO632|Sofia|(SOF)|LOT |E195|(SP-LNG)|3:19 PM|0 hours, -16 minutes
443|442.0|Monday, Jan 09|3:40 PM|LO572|Belgrade|(BEG)|LOT |E195|(SP-LNO)|3:33 PM|0 hours, -7 minutes
444|443.0|Monday, Jan 09|3:45 PM|LO3944|Poznan|(POZ)|LOT |E170|(SP-LDG)|3:22 PM|0 hours, -23 minutes
445|444.0|Monday, Jan 09|3:50 PM|LO3802|Rzeszow|(RZE)|LOT |E75S|(SP-LIA)|3:35 PM|0 hours, -15 minutes
446|445.0|Monday, Jan 09|3:50 PM|LO3838|Gdansk|(GDN)|LOT |E170|(SP-LDI)|3:31 PM|0 hours, -19 minutes
447|446.0|Monday, Jan 09|3:50 PM|LO3844|Wroclaw|(WRO)|LOT |E190|(SP-LME)|3:25 PM|0 hours, -25 minutes
448|447.0|Monday, Jan 09|3:50 PM|LO3906|Krakow|(KRK)|LOT |E190|(SP-LMH)|4:04 PM|0 hours, 14 minutes
449|448.0|Monday, Jan 09|3:50 PM|LO3994|Bydgoszcz|(BZG)|LOT |E170|(SP-LDF)|3:44 PM|0 hours, -6 minutes
450|449.0|Monday, Jan 09|3:58 PM|LPR41|Krakow|(KRK)|Polish Medical Air Rescue |LJ75|(SP-MXR)|5:32 PM|1 hours, 34 minutes
451|450.0|Monday, Jan 09|4:00 PM|LO3884|Katowice|(KTW)|LOT (Sliwka Naleczowska Livery) |E195|(SP-LNC)|4:19 PM|0 hours, 19 minutes
452|451.0|Monday, Jan 09|4:00 PM|LO3934|Szczecin|(SZZ)|LOT |E170|(SP-LDH)|3:46 PM|0 hours, -14 minutes
453|452.0|Monday, Jan 09|4:15 PM|LO264|Amsterdam|(AMS)|LOT |E75S|(SP-LIQ)|3:55 PM|0 hours, -20 minutes
454|453.0|Monday, Jan 09|4:20 PM|W61368|Copenhagen|(CPH)|Wizz Air |A321|(HA-LXK)|4:24 PM|0 hours, 4 minutes
455|454.0|Monday, Jan 09|4:25 PM|AF1046|Paris|(CDG)|Air France |E190|(F-HBLR)|4:42 PM|0 hours, 17 minutes
456|455.0|Monday, Jan 09|4:25 PM|LO438|Barcelona|(BCN)|LOT |B738|(SP-LWA)|3:52 PM|0 hours, -33 minutes
457|456.0|Monday, Jan 09|4:40 PM|CA737|Beijing|(PEK)|Air China |A332|(B-6090)|4:35 PM|0 hours, -5 minutes
458|457.0|Monday, Jan 09|4:45 PM|W61466|Bologna|(BLQ)|Wizz Air |A321|(HA-LXG)|5:29 PM|0 hours, 44 minutes
459|458.0|Monday, Jan 09|4:55 PM|A3872|Athens|(ATH)|Aegean Airlines |A21N|(SX-NAH)|4:48 PM|0 hours, -7 minutes
460|459.0|Monday, Jan 09|5:05 PM|W61372|Eindhoven|(EIN)|Wizz Air |A321|(HA-LXC)|4:59 PM|0 hours, -6 minutes
461|460.0|Monday, Jan 09|5:40 PM|LO3808|Rzeszow|(RZE)|LOT |E75S|(SP-LID)|5:44 PM|0 hours, 4 minutes
462|461.0|Monday, Jan 09|5:45 PM|W61304|London|(LTN)|Wizz Air |A321|(HA-LXN)|5:21 PM|0 hours, -24 minutes
463|462.0|Monday, Jan 09|6:05 PM|LO3924|Krakow|(KRK)|LOT |E190|(SP-LMF)|5:46 PM|0 hours, -19 minutes
"@

$lines = $data -split "`n"
foreach ($line in $lines) {
    $line = $line.TrimEnd("`r")
    if ($line.Trim() -eq "") { continue }
    $parts = $line -split '\|'
    $r = [int]$parts[0]

    $ws.Cells.Item($r, 1).Value = [double]$parts[1]
    $ws.Cells.Item($r, 2).Value = $parts[2]
    $ws.Cells.Item($r, 3).Value = $parts[3]
    $ws.Cells.Item($r, 4).Value = $parts[4]
    $ws.Cells.Item($r, 5).Value = $parts[5]
    $ws.Cells.Item($r, 6).Value = $parts[6]
    $ws.Cells.Item($r, 7).Value = $parts[7]
    $ws.Cells.Item($r, 8).Value = $parts[8]
    $ws.Cells.Item($r, 9).Value = $parts[9]
    $ws.Cells.Item($r, 10).Value = $parts[10]
    $ws.Cells.Item($r, 12).Value = $parts[11]
}

# Columns K (11) and M (13) stay blank for every new row, but still need to
# exist as real (empty) cells, matching the rest of the sheet. Touching the
# border property materializes the cell without allocating a new cell style.
$ws.Range("K425:K463").Borders.LineStyle = -4142
$ws.Range("M425:M463").Borders.LineStyle = -4142
